$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 629.3
$ws.Range("I33").Value = 137.83333
$ws.Range("J33").Value = 1366.5
$ws.Range("K33").Value = 137.83333
$ws.Range("L33").Value = 1366.5
$ws.Range("M33").Value = 91.16667000000001
$ws.Range("N33").Value = -1824.5
$ws.Range("H64").Value = 58630.89
$ws.Range("I64").Value = 127624.75
$ws.Range("J64").Value = 3435.8
$ws.Range("K64").Value = 127624.75
$ws.Range("L64").Value = 3435.8
$ws.Range("M64").Value = -127376.75
$ws.Range("N64").Value = -3931.8
$ws.Range("H67").Value = 58630.89
$ws.Range("I67").Value = 127624.75
$ws.Range("J67").Value = 3435.8
$ws.Range("K67").Value = 127624.75
$ws.Range("L67").Value = 3435.8
$ws.Range("M67").Value = -126766.75
$ws.Range("N67").Value = -5151.8
$ws.Range("H117").Value = 49985
$ws.Range("J117").Value = 49985
$ws.Range("L117").Value = 49985
$ws.Range("N117").Value = -59163
$ws.Range("H138").Value = 3089.6545
$ws.Range("I138").Value = 2611.818
$ws.Range("J138").Value = 3209.1135
$ws.Range("K138").Value = 7835.454000000001
$ws.Range("L138").Value = 9627.3405
$ws.Range("M138").Value = -2695.454000000001
$ws.Range("N138").Value = -19907.3405

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1298.2667
$ws.Range("I45").Value = 924.9091
$ws.Range("J45").Value = 2325
$ws.Range("K45").Value = 924.9091
$ws.Range("L45").Value = 2325
$ws.Range("M45").Value = -547.9091
$ws.Range("N45").Value = -3079
$ws.Range("H53").Value = 9800
$ws.Range("J53").Value = 9800
$ws.Range("L53").Value = 9800
$ws.Range("N53").Value = -11164
$ws.Range("H54").Value = 14000
$ws.Range("J54").Value = 14000
$ws.Range("L54").Value = 14000
$ws.Range("N54").Value = -15538
$ws.Range("H122").Value = 1692.2106
$ws.Range("I122").Value = 1565.1111
$ws.Range("K122").Value = 4695.3333
$ws.Range("M122").Value = -2245.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 20860584
$ws.Range("I107").Value = 27812272
$ws.Range("K107").Value = 27812272
$ws.Range("M107").Value = -27810352
$ws.Range("H116").Value = 35544.332
$ws.Range("J116").Value = 35544.332
$ws.Range("L116").Value = 35544.332
$ws.Range("N116").Value = -44722.332
$ws.Range("H134").Value = 1547.2565
$ws.Range("I134").Value = 1413
$ws.Range("K134").Value = 4239
$ws.Range("M134").Value = -1704

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 43816
$ws.Range("J116").Value = 43816
$ws.Range("L116").Value = 43816
$ws.Range("N116").Value = -52994
$ws.Range("H122").Value = 2922
$ws.Range("I122").Value = 3200.7273
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 9602.1819
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -7152.1819
$ws.Range("N122").Value = -10600
$ws.Range("H134").Value = 1311.3928
$ws.Range("I134").Value = 1013.2917
$ws.Range("J134").Value = 3100
$ws.Range("K134").Value = 3039.8751
$ws.Range("L134").Value = 9300
$ws.Range("M134").Value = -504.8751000000002
$ws.Range("N134").Value = -14370

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2925
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2925
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 8775
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -8943
$ws.Range("H87").Value = 12263.615
$ws.Range("I87").Value = 3271.1667
$ws.Range("K87").Value = 9813.500100000001
$ws.Range("M87").Value = -8565.500100000001
$ws.Range("H90").Value = 12263.615
$ws.Range("I90").Value = 3271.1667
$ws.Range("K90").Value = 29440.5003
$ws.Range("M90").Value = -23200.5003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 35339.332
$ws.Range("J100").Value = 35339.332
$ws.Range("L100").Value = 35339.332
$ws.Range("N100").Value = -37503.332
$ws.Range("H111").Value = 22666.666
$ws.Range("I111").Value = 20000
$ws.Range("J111").Value = 24000
$ws.Range("K111").Value = 20000
$ws.Range("L111").Value = 24000
$ws.Range("M111").Value = -16933
$ws.Range("N111").Value = -30134
$ws.Range("H112").Value = 35989
$ws.Range("J112").Value = 35989
$ws.Range("L112").Value = 35989
$ws.Range("N112").Value = -38205
$ws.Range("H113").Value = 1761.6666
$ws.Range("I113").Value = 1433.6666
$ws.Range("J113").Value = 2089.6667
$ws.Range("K113").Value = 1433.6666
$ws.Range("L113").Value = 2089.6667
$ws.Range("M113").Value = 736.3334
$ws.Range("N113").Value = -6429.6667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 432.64517
$ws.Range("I22").Value = 312.30768
$ws.Range("K22").Value = 312.30768
$ws.Range("M22").Value = -17.30768
$ws.Range("H27").Value = 432.64517
$ws.Range("I27").Value = 312.30768
$ws.Range("K27").Value = 312.30768
$ws.Range("M27").Value = -205.30768
$ws.Range("H48").Value = 13318.333
$ws.Range("J48").Value = 13977.5
$ws.Range("L48").Value = 13977.5
$ws.Range("N48").Value = -15299.5
$ws.Range("H59").Value = 22626.666
$ws.Range("J59").Value = 22626.666
$ws.Range("L59").Value = 22626.666
$ws.Range("N59").Value = -23934.666
$ws.Range("H93").Value = 1882.95
$ws.Range("I93").Value = 1597.0769
$ws.Range("K93").Value = 1597.0769
$ws.Range("M93").Value = -349.0769
$ws.Range("H112").Value = 54599.8
$ws.Range("J112").Value = 54599.8
$ws.Range("L112").Value = 54599.8
$ws.Range("N112").Value = -57553.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 35800
$ws.Range("J57").Value = 35800
$ws.Range("L57").Value = 35800
$ws.Range("N57").Value = -37308
$ws.Range("H74").Value = 10216.75
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 10216.75
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 10216.75
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -12088.75
$ws.Range("H77").Value = 10216.75
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 10216.75
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 30650.25
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -40010.25
$ws.Range("H132").Value = 2914.8438
$ws.Range("I132").Value = 3489.7896
$ws.Range("J132").Value = 2074.5386
$ws.Range("K132").Value = 10469.3688
$ws.Range("L132").Value = 6223.6158
$ws.Range("M132").Value = -7939.3688
$ws.Range("N132").Value = -11283.6158
$ws.Range("H133").Value = 53141.668
$ws.Range("J133").Value = 53141.668
$ws.Range("L133").Value = 53141.668
$ws.Range("N133").Value = -63261.668
